# The workbook is already open.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sampled fault-proneness counts for row 40 ("Hadoop.net"):
#   B40: 255 -> 210
#   C40: 2279 -> 1858
# The TOTAL row (109) contains SUM formulas over B2:B108 / C2:C108, so it
# will automatically recalculate to the new totals (264 / 2272).
$ws.Range("B40").Value = 210
$ws.Range("C40").Value = 1858

# Reflect the author's new cursor position / selection on the sheet
# (previously B109, now C109).
$ws.Range("C109").Select()

$wb.Save()
